# Added Week 15 simulations
$wb = $excel.ActiveWorkbook

# OFF sheet - row 2 updates
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 167
$wsOff.Range("C2").Value = 118
$wsOff.Range("D2").Value = 41
$wsOff.Range("E2").Value = 18
$wsOff.Range("G2").Value = 5

# DEF sheet - row 2 updates
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 199
$wsDef.Range("C2").Value = 138
$wsDef.Range("D2").Value = 50
$wsDef.Range("E2").Value = 26
